# Tweak the London example so the final "plague" rows collapse into a single
# summary row that uses text labels ("plague death toll" / "Second Pandemic")
# instead of the previous numeric mortality-bill figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Write B12 before A12 so the new shared strings are interned in the same
# order as in the target workbook (Second Pandemic, then plague death toll).
$ws.Range("B12").Value = "Second Pandemic"
$ws.Range("A12").Value = "plague death toll"
$ws.Range("C12").Value = 1331
$ws.Range("D12").Value = 1750
# E12 (Level = 1) and F12 (add = FALSE) are unchanged.

# The old rows 13-18 (other mortality-bill entries) are removed entirely,
# shrinking the sheet from A1:F18 down to A1:F12.
$ws.Range("A13:F18").Delete() | Out-Null

# Reflect the new selection left behind in the saved file.
$ws.Range("A12").Select() | Out-Null
